$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 320
$ws1.Range("F6").Value = 399
$ws1.Range("F7").Value = 865
$ws1.Range("F8").Value = 51
$ws1.Range("F9").Value = 516
$ws1.Range("F11").Value = 296
$ws1.Range("F12").Value = 1118
$ws1.Range("F17").Value = 6639
$ws1.Range("F21").Value = 7584
$ws1.Range("I21").Value = "//i0.hdslb.com/bfs/openplatform/202402/WzjKw3B41708484793603.jpeg"
$ws1.Range("F24").Value = 3397
$ws1.Range("F26").Value = 1801
$ws1.Range("F28").Value = 4516
$ws1.Range("F29").Value = 135
$ws1.Range("F30").Value = 352
$ws1.Range("F32").Value = 222
$ws1.Range("F34").Value = 1683
$ws1.Range("F36").Value = 168
$ws1.Range("F39").Value = 1206
$ws1.Range("F40").Value = 1775
$ws1.Range("F41").Value = 2136

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 1226

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1226
$ws4.Range("F7").Value = 320
$ws4.Range("F8").Value = 399
$ws4.Range("F9").Value = 865
$ws4.Range("F10").Value = 51
$ws4.Range("F11").Value = 516
$ws4.Range("F13").Value = 296
$ws4.Range("F14").Value = 1119
$ws4.Range("F20").Value = 6639
$ws4.Range("F24").Value = 7584
$ws4.Range("I24").Value = "//i0.hdslb.com/bfs/openplatform/202402/WzjKw3B41708484793603.jpeg"
$ws4.Range("F27").Value = 3397
$ws4.Range("F29").Value = 1801
$ws4.Range("F31").Value = 4516
$ws4.Range("F32").Value = 135
$ws4.Range("F33").Value = 352
$ws4.Range("F36").Value = 222
$ws4.Range("F38").Value = 1683
$ws4.Range("F40").Value = 168
$ws4.Range("F44").Value = 1206
$ws4.Range("F45").Value = 1775
$ws4.Range("F47").Value = 2136
